# Apply the OOXML diff: insert two new index columns (D,E) before the
# existing "Unnamed: 0.1" column, duplicate the row index into the two
# freshly-vacated columns (G,H), and special-case the "Meadowlark" row
# (16) whose venue match was dropped in the refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank columns at D:E -- this shifts the old D..P block to
#    F..R (matches how the header row names moved in the diff).
$ws.Columns("D:E").Insert()

# 2. New header cells created by the insert need their own text; Excel
#    already carried "Unnamed: 0.1".."location_geom" two slots to the
#    right, so we only need to fill in the two brand-new columns.
$ws.Range("D1").Value = "Unnamed: 0.3"
$ws.Range("E1").Value = "Unnamed: 0.2"

# 3. For every data row, columns G and H are brand-new (blank after the
#    insert). In the refreshed workbook they simply repeat the same row
#    index that already sits in D/E/F -- except for rows whose venue
#    match is entirely missing, where every one of D..R is blank.
$lastRow = 22
$blankRows = @(9, 11, 21)

for ($r = 2; $r -le $lastRow; $r++) {
    if ($blankRows -contains $r) {
        continue
    }

    $idx = $ws.Range("F$r").Value2
    if ($idx -ne $null) {
        $ws.Range("G$r").Value = $idx
        $ws.Range("H$r").Value = $idx
    }
}

# 4. Row 16 ("Meadowlark") no longer has a matched venue record in the
#    refreshed source: its supersite name now carries the old gmap_name,
#    and every venue-derived column D..R is cleared out.
$ws.Range("B16").Value = "Meadowlark School"
$ws.Range("D16:R16").ClearContents()

# 5. Keep the sheet's declared dimension in sync with the new R-column
#    extent.
$ws.Range("A1:R22").Select()
